$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: append section-4 facility data columns (AX10:BR10)
$ws.Range("AX10").Value = 'We have an emergency action plan written up. Speed dial set to 911, AED system for facility, walker talkies on staff leaders.'
$ws.Range("AY10").Value = 'No'
$ws.Range("AZ10").Value = 'Yes'
$ws.Range("BA10").Value = 'Yes'
$ws.Range("BB10").Value = 'Gravel/dirt road,Forestry road,Private access/limited access road,Extra large vehicle (Such semi-tuck with trailer, standard bus, class A RV),Large vehicles (such as semi without trailer, short bus length, Class C RV, truck with trailer),Standard vehicle (Such as truck without trailer, SUV, sedan),Small vehicle (Such as motorcycle, hybrid, smart car)'
$ws.Range("BD10").Value = '104 acres'
$ws.Range("BE10").Value = 'Stream, river, marsh, forest, reptile and amphibians habitat, wildlife, gold beds, fossil beds'
$ws.Range("BF10").Value = 'Food service facility, kitchen staff, and food'
$ws.Range("BH10").Value = 'Yes (if yes, please check all that apply),Vegetarian,Vegan,Kosher,Diabetic Meal plan,Gluten/wheat allergies,Dairy allergies (milk and/or egg),Soy, peanut, or tree nut allergies,Seafood or shellfish allergies'
$ws.Range("BJ10").Value = 240
$ws.Range("BK10").Value = 300
$ws.Range("BL10").Value = 'If Yes, How many and what capacity?'
$ws.Range("BM10").Value = 5
$ws.Range("BN10").Value = 180
$ws.Range("BO10").Value = 500
$ws.Range("BP10").Value = 'Yes, we have sidewalks to all buildings with 2 ADA accessible cabins, and ADA accessible batthroom facilities.'
$ws.Range("BQ10").Value = 'More bathroom facilities.'
$ws.Range("BR10").Value = 'We are currently providing curriculum and turn key outdoor school for schools with limited resources to plan their own.'

# Row 12: replace with correct section 1-3 response data (A12:T12); clear stray K12/M12
$ws.Range("A12").Value = 43381.7379050926
$ws.Range("B12").Value = 43381.75645833334
$ws.Range("C12").Value = 'IP Address'
$ws.Range("D12").Value = '172.58.35.118'
$ws.Range("E12").Value = 100
$ws.Range("F12").Value = 1603
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = 43381.75645833334
$ws.Range("I12").Value = 'R_21yUQgUKxVX6Z8y'
$ws.Range("N12").Value = 38.4797058105468
$ws.Range("O12").Value = -121.443801879882
$ws.Range("P12").Value = 'anonymous'
$ws.Range("Q12").Value = 'EN'
$ws.Range("R12").Value = 'joel@grovechristiancamp.org'
$ws.Range("S12").Value = 'Grove Christian Camp / Grove Camp Outdoor School'
$ws.Range("T12").Value = 1
$ws.Range("K12").Value = ""
$ws.Range("M12").Value = ""

# Row 13: append section-4 facility data columns (AX13:BP13)
$ws.Range("AX13").Value = 'We are located in the Siusalaw National Forest and are "off the grid". However, we do have a satellite communication system that provides both telephone service and internet access. The local hospital, Pac West and Life Flight all have our coordinates and are familiar with our location in case of an emergency. '
$ws.Range("AY13").Value = 'No'
$ws.Range("AZ13").Value = 'Yes'
$ws.Range("BA13").Value = 'No'
$ws.Range("BB13").Value = 'Paved road,Gravel/dirt road,Extra large vehicle (Such semi-tuck with trailer, standard bus, class A RV),Large vehicles (such as semi without trailer, short bus length, Class C RV, truck with trailer),Standard vehicle (Such as truck without trailer, SUV, sedan),Small vehicle (Such as motorcycle, hybrid, smart car),Horse/Pack animal'
$ws.Range("BD13").Value = '25 acres'
$ws.Range("BE13").Value = 'Forest (old growth forest including Sitka Spruce and many nurse logs); pristine mountain stream (including salmon spawning location); many ferns and wild plants and countless mushrooms.  '
$ws.Range("BF13").Value = 'Other'
$ws.Range("BG13").Value = 'Programs can choose to provide their own food service or hire DCC''s staff to provide meals.'
$ws.Range("BH13").Value = 'Yes (if yes, please check all that apply),Vegetarian,Gluten/wheat allergies,Dairy allergies (milk and/or egg),Soy, peanut, or tree nut allergies'
$ws.Range("BJ13").Value = 150
$ws.Range("BK13").Value = 'Depends on the weather.'
$ws.Range("BL13").Value = 'If Yes, How many and what capacity?'
$ws.Range("BM13").Value = 'We have a meeting area that can seat 100 and a craft room that can seat about 12. There is also an outdoor covered Activity Center that could seat 100+. '
$ws.Range("BN13").Value = 'The main lodge can sleep up to 120. This would include upper bunks. Rooms are chalet-style with common restrooms at each end of the hallway. In addition to the 3-floor lodge, we also have 5 cabins that each sleep 12-14 in bunk beds and a yurt that sleeps 14. A central bathhouse is available for cabin campers. The total capacity is about 190. '
$ws.Range("BP13").Value = 'Yes, we have two main floor rooms available in the lodge that are accessible and the Alsea cabin is also accessible. '

# Row 15: replace with correct section 1-3 response data (A15:T15); clear stray J15/K15/M15
$ws.Range("A15").Value = 43382.42454861111
$ws.Range("B15").Value = 43382.44295138889
$ws.Range("C15").Value = 'IP Address'
$ws.Range("D15").Value = '69.168.127.185'
$ws.Range("E15").Value = 100
$ws.Range("F15").Value = 1590
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 43382.44296296296
$ws.Range("I15").Value = 'R_1rMqBg85bkdvfln'
$ws.Range("N15").Value = 45.0897064208984
$ws.Range("O15").Value = -123.400299072265
$ws.Range("P15").Value = 'anonymous'
$ws.Range("Q15").Value = 'EN'
$ws.Range("R15").Value = 'info@driftcreek.org'
$ws.Range("S15").Value = 'Drift Creek Camp'
$ws.Range("T15").Value = 1
$ws.Range("J15").Value = ""
$ws.Range("K15").Value = ""
$ws.Range("M15").Value = ""
